$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new row of data (2nd Sunday of Advent - mode 1, differentia "g")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "g"
$ws.Range("C2").Value = "c4"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "(hr g) " + [char]0x2020 + "(;3) (hr ixi g h) <v>`$\star`$</v>(;) (hr g f g ::)"

# Update the active selection to F2
$ws.Range("F2").Select()
